
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ランサーズ")
$ws2 = $wb.Worksheets.Item("統計")

# ---- Sheet1 (ランサーズ): insert 3 new rows at the top of the data (rows 2-4) ----
$ws1.Range("A2:A4").EntireRow.Insert()

# Fill in the 3 new rows with the new scraped entries
$ws1.Range("A2").Value2 = "2025-08-31 01:21:48"
$ws1.Range("B2").Value2 = "【RPA自動化】通関データ入力業務の効率化依頼"
$ws1.Range("C2").Value2 = "システム開発"
$ws1.Range("D2").Value2 = "50,000 円 ~ 100,000 円 / 固定"
$ws1.Range("E2").Value2 = "期限情報なし"
$ws1.Range("F2").Value2 = "https://www.lancers.jp/work/detail/5383211"
$ws1.Range("G2").Value2 = 153
$ws1.Range("H2").Value2 = "◆効率化,自動化"

$ws1.Range("A3").Value2 = "2025-08-31 01:21:48"
$ws1.Range("B3").Value2 = "IB報酬を得るための高性能EA開発依頼"
$ws1.Range("C3").Value2 = "システム開発"
$ws1.Range("D3").Value2 = "100,000 円 ~ 200,000 円 / 固定"
$ws1.Range("E3").Value2 = "期限情報なし"
$ws1.Range("F3").Value2 = "https://www.lancers.jp/work/detail/5383199"
$ws1.Range("G3").Value2 = 68
$ws1.Range("H3").Value2 = "◆開発"

$ws1.Range("A4").Value2 = "2025-08-31 01:21:48"
$ws1.Range("B4").Value2 = "マクロの仕様変更をお願いします。"
$ws1.Range("C4").Value2 = "システム開発"
$ws1.Range("D4").Value2 = "~ 5,000 円 / 固定"
$ws1.Range("E4").Value2 = "期限情報なし"
$ws1.Range("F4").Value2 = "https://www.lancers.jp/work/detail/5383131"
$ws1.Range("G4").Value2 = 10

# ---- Rebuild hyperlinks on column F so refs/ids line up with the shifted rows ----
$ws1.Hyperlinks.Delete()
for ($r = 2; $r -le 70; $r++) {
  $cell = $ws1.Range("F" + $r)
  $url = $cell.Value2
  if ($url -ne $null -and $url -ne "") {
    $ws1.Hyperlinks.Add($cell, $url) | Out-Null
  }
}

# ---- Sheet2 (統計): append new summary row ----
$ws2.Range("A31").Value2 = "2025-08-31T01:21:48.163127"
$ws2.Range("B31").Value2 = 13
$ws2.Range("C31").Value2 = "全案件リスト"
$ws2.Range("D31").Value2 = 76.90000000000001
$ws2.Range("E31").Value2 = 3
$ws2.Range("F31").Value2 = 7
$ws2.Range("G31").Value2 = 13

Write-Host "done"
